$d = $word.ActiveDocument

# 1. Re-set the run text for "Mara (stretching satisfaction):" so Word
#    regenerates the <w:t> element (picking up xml:space="preserve").
$d.Content.Find.Execute("Mara (stretching satisfaction):", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Mara (stretching satisfaction):", 2) | Out-Null
